# Refresh the cryptocurrency price/volume snapshot (columns D and E,
# rows 2-51) as produced by the scheduled GitHub Actions data pull.
#
# Column D values that look like plain decimal numbers are entered with a
# leading apostrophe so Excel keeps them as literal text (matching the
# source data, e.g. "1.000" must stay "1.000" and not collapse to "1").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.387.20"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "1.848.31"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'240.50"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.07636"
$ws.Range("E8").Value = "  +0.47%  "
$ws.Range("E9").Value = "  -0.93%  "
$ws.Range("D10").Value = "'24.71"
$ws.Range("E10").Value = "  +0.94%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "'5.024"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Value = "'0.6788"
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").Value = "'0.00001064"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").Value = "'83.18"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "29.404.52"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D18").Value = "'226.72"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'7.501"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'158.12"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'0.1378"
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  +5.70%  "
$ws.Range("D28").Value = "'1.465"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("D29").Value = "'0.05597"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").Value = "'4.120"
$ws.Range("E30").Value = "  +0.27%  "
$ws.Range("D31").Value = "'4.078"
$ws.Range("E31").Value = "  +0.89%  "
$ws.Range("D32").Value = "'1.838"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").Value = "'0.6946"
$ws.Range("E34").Value = "  -2.14%  "
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "1.231.12"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "'0.01801"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'2.716"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "'6.397"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("D40").Value = "'0.9048"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("D42").Value = "'101.58"
$ws.Range("D43").Value = "'66.06"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("D46").Value = "'0.4014"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").Value = "'9.005"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").Value = "'1.680"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'0.1142"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("D50").Value = "'0.05702"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  +0.03%  "
